$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before assigning, to preserve
# the exact text representation (e.g. trailing zeros like "1.000").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.712.28'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.889.08'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '249.52'
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '0.4759'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.2935'
$ws.Range("E8").Value = '  +0.69%  '
$ws.Range("D9").Value = '0.06535'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '22.02'
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("D11").Value = '0.07747'
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("D12").Value = '0.7411'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '96.90'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").Value = '1.887.86'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").Value = '5.257'
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").Value = '276.09'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").Value = '30.792.60'
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("D18").Value = '13.22'
$ws.Range("E18").Value = '  -2.89%  '
$ws.Range("D19").Value = '0.000007565'
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = '2.134.39'
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").Value = '5.328'
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("D23").Value = '0.9996'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '6.242'
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").Value = '9.244'
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("D26").Value = '163.99'
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("D27").Value = '18.83'
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("D28").Value = '1.929'
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("D29").Value = '1.345'
$ws.Range("E29").Value = '  -2.28%  '
$ws.Range("D30").Value = '0.09740'
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("D31").Value = '1.507'
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("D32").Value = '4.308'
$ws.Range("E32").Value = '  -0.61%  '
$ws.Range("D33").Value = '4.199'
$ws.Range("E33").Value = '  +2.68%  '
$ws.Range("D34").Value = '0.04884'
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("D35").Value = '1.128'
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D36").Value = '0.7012'
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("D37").Value = '2.722'
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Value = '0.01922'
$ws.Range("E38").Value = '  +2.48%  '
$ws.Range("D39").Value = '2.802'
$ws.Range("E39").Value = '  +2.63%  '
$ws.Range("D40").Value = '6.347'
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").Value = '76.03'
$ws.Range("E41").Value = '  +6.22%  '
$ws.Range("D42").Value = '2.033'
$ws.Range("E42").Value = '  +3.46%  '
$ws.Range("D43").Value = '0.4260'
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("D44").Value = '0.8417'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '0.9998'
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '102.35'
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("D47").Value = '9.396'
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("D48").Value = '7.074'
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").Value = '35.72'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = '922.48'
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("D51").Value = '0.05775'
$ws.Range("E51").Value = '  +2.31%  '

# Reset the style on column D back to Normal so no stray formatting
# (other than the original General) is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
